# Revert "Changes made for demo":
# Remove the extra demo row (row 17) that was added to the "Excluded
# structures" sheet - including the "for testing" placeholder text that
# was entered in columns A-G, J & K, and the sample category picks that
# were entered in columns H & I.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Excluded structures")
$ws.Activate()

# Columns A-G keep their wrapped-text styling on every data row, so only
# blank out the values (leave the formatting in place), matching rows 2-16.
$ws.Range("A17:G17").ClearContents()

# Columns H-K never carried any formatting of their own on row 17 - fully
# clear them (contents + formats) so the cells disappear entirely, same as
# an empty/untouched row.
$ws.Range("H17:K17").Clear()

# Restore the selection to where it was before the demo edits (column I,
# row 17, inside the frozen header pane) instead of column K.
$ws.Range("I17").Select() | Out-Null
